$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.587.80"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.600.04"
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.34"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.90"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +2.14%  "

$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("E11").Value = "  +0.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.85"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  +1.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.597.00"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  +0.94%  "

$ws.Range("E14").Value = "  +3.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.606.56"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.53"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  +2.17%  "

$ws.Range("E19").Value = "  +2.48%  "

$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("E23").Value = "  +0.68%  "

$ws.Range("E24").Value = "  -0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.02"
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = "  +1.43%  "

$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("E28").Value = "  +0.72%  "

$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("E30").Value = "  +1.95%  "

$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  +0.36%  "

$ws.Range("E33").Value = "  +2.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.421.31"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  +1.97%  "

$ws.Range("E36").Value = "  +3.84%  "

$ws.Range("E37").Value = "  -1.90%  "

$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("E39").Value = "  +3.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.543"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  +2.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.04"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  +6.01%  "

$ws.Range("E42").Value = "  +0.27%  "

$ws.Range("E43").Value = "  +5.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.809"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  +2.47%  "

$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.989"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  +16.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "66.11"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  +2.11%  "

$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.739.44"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  +1.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.19"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  +0.68%  "

$ws.Range("E51").Value = "  +0.89%  "
